# LM6-Week-12.docx restyle:
#  - Apply "Inter" font to every run in the document (body + table).
#  - Apply heading color #0F4761 to every bold run (body + table headers).
#  - Fix "Each classmate response is worth 2 points" -> "... worth 10 points".
#  - Shrink the rubric table to ~97.3% (tblW 5000->4865 pct, columns 2640->2568 dxa).
#  - Set all page margins to 0.5" (720 twips).

$d = $word.ActiveDocument

$COLOR_0F4761 = 6375183   # BGR long for RGB 0F4761 (Word Font.Color expects BGR)

# ---------------------------------------------------------------------------
# 1) Font: set "Inter" across the whole body range.
# ---------------------------------------------------------------------------
$endPos = $d.Content.End
$bodyRng = $d.Range(0, $endPos)
$bodyRng.Font.Name = "Inter"

# Range.Font does not recurse into table content on its own, so handle tables
# explicitly too.
foreach ($tbl in $d.Tables) {
    $tbl.Range.Font.Name = "Inter"
}

# ---------------------------------------------------------------------------
# 2) Color: every bold run gets #0F4761 (heading / emphasis color).
# ---------------------------------------------------------------------------
$find = $bodyRng.Find
$find.ClearFormatting()
$find.Font.Bold = $true
$find.Text = ""
$find.Forward = $true
$find.Wrap = 0
$i = 0
while ($find.Execute() -and $i -lt 100) {
    $i = $i + 1
    if ($bodyRng.Start -eq $bodyRng.End) {
        break
    }
    $bodyRng.Font.Color = $COLOR_0F4761
    $bodyRng.Collapse(0)
}

# Table header cells ("Exceptional (2 pts)", "Emerging (1 pt)", "Missing (0 pts)")
# are bold too, but Find above does not reach into table cells, so color them
# directly.
foreach ($tbl in $d.Tables) {
    foreach ($row in $tbl.Rows) {
        foreach ($cell in $row.Cells) {
            $cellRng = $cell.Range
            if ($cellRng.Font.Bold -eq -1) {
                $cellRng.Font.Color = $COLOR_0F4761
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 3) Text fix: "worth 2 points" -> "worth 10 points" in the assessment blurb.
# ---------------------------------------------------------------------------
$fixRng = $d.Content
$fixRng.Find.ClearFormatting()
$fixRng.Find.Replacement.ClearFormatting()
$fixRng.Find.Execute(
    "Each classmate response is worth 2 points",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Each classmate response is worth 10 points", 2
)

# ---------------------------------------------------------------------------
# 4) Rubric table: shrink to ~97.3% (5000 -> 4865 pct / 2640 -> 2568 dxa).
# ---------------------------------------------------------------------------
$tbl = $d.Tables.Item(1)
$tbl.PreferredWidth = 243.25   # engine units = (w:w twentieths-of-percent) / 20 -> 4865/20
foreach ($col in $tbl.Columns) {
    $col.Width = 128.4         # points -> 128.4 * 20 = 2568 dxa
}

# ---------------------------------------------------------------------------
# 5) Page margins: 0.5" (720 twips = 36 pt) on all sides.
# ---------------------------------------------------------------------------
foreach ($sec in $d.Sections) {
    $sec.PageSetup.TopMargin = 36
    $sec.PageSetup.BottomMargin = 36
    $sec.PageSetup.LeftMargin = 36
    $sec.PageSetup.RightMargin = 36
}
